$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column A stays text (so the dashed dates are not auto-converted
# into date serial numbers by Excel's input parser).
$ws.Range("A3:A21").NumberFormat = "@"

# Update dates (slash -> dash) and attendance flag values for each row
# Row 3
$ws.Range("A3").Value = "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4
$ws.Range("A4").Value = "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 5
$ws.Range("A5").Value = "04-08-2022"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# Row 6
$ws.Range("A6").Value = "08-08-2022"

# Row 7
$ws.Range("A7").Value = "11-08-2022"

# Row 8
$ws.Range("A8").Value = "15-08-2022"

# Row 9
$ws.Range("A9").Value = "18-08-2022"

# Row 10
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("H10").Value = 0

# Row 11
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("H11").Value = 0

# Row 12
$ws.Range("A12").Value = "29-08-2022"

# Row 13
$ws.Range("A13").Value = "01-09-2022"

# Row 14
$ws.Range("A14").Value = "05-09-2022"

# Row 15
$ws.Range("A15").Value = "08-09-2022"

# Row 16
$ws.Range("A16").Value = "12-09-2022"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("H16").Value = 0

# Row 17
$ws.Range("A17").Value = "15-09-2022"

# Row 18
$ws.Range("A18").Value = "19-09-2022"

# Row 19
$ws.Range("A19").Value = "22-09-2022"

# Row 20
$ws.Range("A20").Value = "26-09-2022"

# Row 21
$ws.Range("A21").Value = "29-09-2022"
